$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.03208162110170181
$ws.Range("H2").Value = -3.644261968430065
$ws.Range("I2").Value = -8.625475708809777
$ws.Range("G3").Value = 0.06687162342617554
$ws.Range("H3").Value = 32.51168879827196
$ws.Range("G4").Value = 0.0139640214175768
$ws.Range("H4").Value = -71.31220213840888
$ws.Range("G5").Value = 0.07123469294217934
$ws.Range("H5").Value = 0.8211045158787089
$ws.Range("G6").Value = -0.1540462429413415
$ws.Range("H6").Value = -30.96310630252292
$ws.Range("G7").Value = -0.07592735945351084
$ws.Range("H7").Value = 39.27437248891227
$ws.Range("G8").Value = -0.2188718529334285
$ws.Range("H8").Value = -9.88716516511186
$ws.Range("G9").Value = -0.3033111385678904
$ws.Range("H9").Value = -0.4085884864675501
$ws.Range("G10").Value = -0.01685330631372322
$ws.Range("H10").Value = -1424.201084290179
$ws.Range("G11").Value = -0.01782258904769769
$ws.Range("H11").Value = 13.32852311274608
$ws.Range("G12").Value = 0.1997671648356183
$ws.Range("H12").Value = -5.73075344194351
$ws.Range("G13").Value = 0.2466165245753381
$ws.Range("H13").Value = 4.944725032814968
$ws.Range("G14").Value = -0.07377592133293376
$ws.Range("H14").Value = 18.97231949575492
$ws.Range("G15").Value = -0.05251792841898704
$ws.Range("H15").Value = 25.99849703899318
$ws.Range("G16").Value = 0.1596809718707339
$ws.Range("H16").Value = -16.5953318851737
$ws.Range("G17").Value = 0.1719701244552567
$ws.Range("H17").Value = -1.092484591019351
$ws.Range("G18").Value = 0.0555653284199867
$ws.Range("H18").Value = 2.462614242737266
$ws.Range("G19").Value = 0.06316901527490057
$ws.Range("H19").Value = -26.58146443679704
$ws.Range("G20").Value = 0.001705204562111364
$ws.Range("H20").Value = -86.60368507898536
$ws.Range("G21").Value = -0.02984487959186669
$ws.Range("H21").Value = 44.6022268960774
$ws.Range("G22").Value = 0.07881842909534316
$ws.Range("H22").Value = 20.75244432245246
$ws.Range("G23").Value = 0.07171703472764612
$ws.Range("H23").Value = 24.35030426984009
$ws.Range("G24").Value = 0.05927091326668833
$ws.Range("H24").Value = 82.96817773240419
$ws.Range("G25").Value = 0.03725331983404287
$ws.Range("H25").Value = 26.57364362852677
$ws.Range("G26").Value = 0.115548077803632
$ws.Range("H26").Value = 1.99191042628942
$ws.Range("G27").Value = 0.1163247761972346
$ws.Range("H27").Value = 28.98081563101872
$ws.Range("G28").Value = 0.0907411455445794
$ws.Range("H28").Value = -22.75187476016962
$ws.Range("G29").Value = 0.1248529527507469
$ws.Range("H29").Value = 4.367959754699173
$ws.Range("G30").Value = 0.07283114963935999
$ws.Range("H30").Value = 8.334010045006417
$ws.Range("G31").Value = 0.06063925763455246
$ws.Range("H31").Value = -11.62733353295091
$ws.Range("G32").Value = 0.03795619933500401
$ws.Range("H32").Value = -13.07718572990389
$ws.Range("G33").Value = 0.07224072383406741
$ws.Range("H33").Value = 32.94596508461807
$ws.Range("G34").Value = -0.02694017671003928
$ws.Range("H34").Value = -41.05824902653418
$ws.Range("G35").Value = 0.02393126726749577
$ws.Range("H35").Value = 71.51749255942964
$ws.Range("G36").Value = -0.002577753549880305
$ws.Range("H36").Value = -116.6748645802464
$ws.Range("G37").Value = -0.006956544824441795
$ws.Range("H37").Value = -155.5471679600986
$ws.Range("G38").Value = 0.07551317751573167
$ws.Range("H38").Value = 5.255409359802937
$ws.Range("G39").Value = 0.03308809686384414
$ws.Range("H39").Value = -23.16263037591851
$ws.Range("G40").Value = 0.05757831748328093
$ws.Range("H40").Value = 28.72483355685662
$ws.Range("G41").Value = 0.06844463873485879
$ws.Range("H41").Value = 453.7895800918752
$ws.Range("G42").Value = 0.07930459735437771
$ws.Range("H42").Value = 51.69500217900982
$ws.Range("G43").Value = 0.06216074181148436
$ws.Range("H43").Value = 24.57091124021361
$ws.Range("G44").Value = 0.133801338737148
$ws.Range("H44").Value = 1.556325815849139
$ws.Range("G45").Value = 0.1649778194509112
$ws.Range("H45").Value = -8.052797001414502
$ws.Range("G46").Value = -0.02581450652603963
$ws.Range("H46").Value = 41.23696937017581
$ws.Range("G47").Value = 0.007244012978369976
$ws.Range("H47").Value = 376.5347692307329
$ws.Range("G48").Value = 0.01534174556852695
$ws.Range("H48").Value = 5.864858727421331
$ws.Range("G49").Value = -0.005370105665969714
$ws.Range("H49").Value = 3.394836036202156
$ws.Range("G50").Value = 0.1360471527627251
$ws.Range("H50").Value = -4.816516499381359
$ws.Range("G51").Value = 0.1272392876536723
$ws.Range("H51").Value = -2.84497448326569
$ws.Range("G52").Value = 0.08948338228483049
$ws.Range("H52").Value = 44.43821827429733
$ws.Range("G53").Value = 0.09398543198388534
$ws.Range("H53").Value = 53.61708153716441
$ws.Range("G54").Value = -0.1473785421461656
$ws.Range("H54").Value = -65.02473647869907
$ws.Range("G55").Value = -0.0385111187600477
$ws.Range("H55").Value = 62.89394306092986
$ws.Range("G56").Value = 0.1535480776234136
$ws.Range("H56").Value = -1.015671903473168
$ws.Range("G57").Value = 0.167907724562929
$ws.Range("H57").Value = 20.37595686674489
